$wb = $excel.ActiveWorkbook

$wsHierarchy = $wb.Worksheets.Item("Product_Hierarchy")
$wsExample   = $wb.Worksheets.Item("Example")

# "Product Line Group Code" header is being renamed to
# "Existing Product Line Group Code" on both sheets (A1).
$wsHierarchy.Range("A1").Value = "Existing Product Line Group Code"
$wsExample.Range("A1").Value = "Existing Product Line Group Code"

# Product_Hierarchy!D2 / D3 used to both read "E1&SAP" -- split them
# into the two distinct values now used ("SAP & E1" / "SAP & E2").
$wsHierarchy.Range("D2").Value = "SAP & E1"
$wsHierarchy.Range("D3").Value = "SAP & E2"

# Restore the selections recorded in each sheet view.
$wsExample.Activate() | Out-Null
$wsExample.Range("B2").Select() | Out-Null

$wsHierarchy.Activate() | Out-Null
$wsHierarchy.Range("D2:D3").Select() | Out-Null
